$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 53: Ford Model T intro year correction 1908 -> 1912 (formula recalculates automatically)
$ws.Range("B53").Value = 1912

# Row 54: Mercedes-Benz Citaro (new bus)
$ws.Range("A54").Value = "Mercedes-Benz Citaro"
$ws.Range("B54").Value = 2001
$ws.Range("C54").Value = 2
$ws.Range("D54").Value = "Bus"
$ws.Range("E54").Formula = "=IF(B54 > 1900, ((B54-1900)*10)+400+C54, ((B54-1730)*2)+C54)+VLOOKUP(D54,'ID Scheme'!`$A`$2:`$B`$6,2, FALSE)"
$ws.Range("F54").Value = 55
$ws.Range("G54").Value = 140
$ws.Range("H54").Formula = "=SQRT(F54*G54)/`$B`$1"
$ws.Range("I54").Formula = "=H54*0.9"
$ws.Range("J54").Value = "x"
$ws.Range("H54:J54").NumberFormat = "0"

# Row 55: Leyland Leopard (new bus)
$ws.Range("A55").Value = "Leyland Leopard"
$ws.Range("B55").Value = 1958
$ws.Range("C55").Value = 2
$ws.Range("D55").Value = "Bus"
$ws.Range("E55").Formula = "=IF(B55 > 1900, ((B55-1900)*10)+400+C55, ((B55-1730)*2)+C55)+VLOOKUP(D55,'ID Scheme'!`$A`$2:`$B`$6,2, FALSE)"
$ws.Range("F55").Value = 55
$ws.Range("G55").Value = 57
$ws.Range("H55").Formula = "=SQRT(F55*G55)/`$B`$1"
$ws.Range("I55").Formula = "=H55*0.9"
$ws.Range("J55").Value = "x"
$ws.Range("H55:J55").NumberFormat = "0"

# Update selection to match target view
$ws.Range("G54").Select()
